# Auto-generated edit script: updates recalculated price/profit values
# in the Masamune_Profits sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H109").Value = 26733.6
$ws.Range("J109").Value = 26733.6
$ws.Range("L109").Value = 26733.6
$ws.Range("N109").Value = -29507.6

$ws.Range("H114").Value = 35801.332
$ws.Range("J114").Value = 35801.332
$ws.Range("L114").Value = 35801.332
$ws.Range("N114").Value = -44479.332

$ws.Range("H117").Value = 34060
$ws.Range("J117").Value = 34060
$ws.Range("L117").Value = 34060
$ws.Range("N117").Value = -43238

$ws.Range("H124").Value = 48511.668
$ws.Range("J124").Value = 48511.668
$ws.Range("L124").Value = 48511.668
$ws.Range("N124").Value = -58331.668

$ws.Range("H128").Value = 44170.5
$ws.Range("J128").Value = 44170.5
$ws.Range("L128").Value = 44170.5
$ws.Range("N128").Value = -54130.5

$ws.Range("H130").Value = 45265.6
$ws.Range("J130").Value = 45265.6
$ws.Range("L130").Value = 45265.6
$ws.Range("N130").Value = -55305.6

$ws.Range("H137").Value = 3610.7646
$ws.Range("I137").Value = 1001.1539
$ws.Range("J137").Value = 6324.76
$ws.Range("K137").Value = 3003.4617
$ws.Range("L137").Value = 18974.28
$ws.Range("M137").Value = -453.4616999999998
$ws.Range("N137").Value = -24074.28

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H113").Value = 43000
$ws.Range("J113").Value = 43000
$ws.Range("L113").Value = 43000
$ws.Range("N113").Value = -51678

$ws.Range("H117").Value = 43826.43
$ws.Range("J117").Value = 43826.43
$ws.Range("L117").Value = 43826.43
$ws.Range("N117").Value = -53004.43

$ws.Range("H118").Value = 36909
$ws.Range("J118").Value = 36909
$ws.Range("L118").Value = 36909
$ws.Range("N118").Value = -40223

$ws.Range("H123").Value = 35614.5
$ws.Range("J123").Value = 35614.5
$ws.Range("L123").Value = 35614.5
$ws.Range("N123").Value = -45414.5

$ws.Range("H125").Value = 49297.332
$ws.Range("J125").Value = 49297.332
$ws.Range("L125").Value = 49297.332
$ws.Range("N125").Value = -59137.332

$ws.Range("H130").Value = 39830.668
$ws.Range("J130").Value = 39830.668
$ws.Range("L130").Value = 39830.668
$ws.Range("N130").Value = -49870.668

$ws.Range("H131").Value = 51613
$ws.Range("J131").Value = 51613
$ws.Range("L131").Value = 51613
$ws.Range("N131").Value = -61693

$ws.Range("H137").Value = 36657.8
$ws.Range("J137").Value = 45645
$ws.Range("L137").Value = 45645
$ws.Range("N137").Value = -55845

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H57").Value = 69999
$ws.Range("J57").Value = 69999
$ws.Range("L57").Value = 69999
$ws.Range("N57").Value = -71439

$ws.Range("H117").Value = 41998.223
$ws.Range("J117").Value = 41998.223
$ws.Range("L117").Value = 41998.223
$ws.Range("N117").Value = -51176.223

$ws.Range("H124").Value = 47996
$ws.Range("J124").Value = 47996
$ws.Range("L124").Value = 47996
$ws.Range("N124").Value = -57816

$ws.Range("H125").Value = 50780
$ws.Range("J125").Value = 50780
$ws.Range("L125").Value = 50780
$ws.Range("N125").Value = -60620

$ws.Range("H126").Value = 42241.332
$ws.Range("J126").Value = 42241.332
$ws.Range("L126").Value = 42241.332
$ws.Range("N126").Value = -52121.332

$ws.Range("H130").Value = 48275.75
$ws.Range("J130").Value = 48275.75
$ws.Range("L130").Value = 48275.75
$ws.Range("N130").Value = -58315.75

$ws.Range("H136").Value = 69999
$ws.Range("J136").Value = 69999
$ws.Range("L136").Value = 69999
$ws.Range("N136").Value = -80199

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H20").Value = 49319.6
$ws.Range("J20").Value = 49319.6
$ws.Range("L20").Value = 49319.6
$ws.Range("N20").Value = -49791.6

$ws.Range("H30").Value = 49319.6
$ws.Range("J30").Value = 49319.6
$ws.Range("L30").Value = 49319.6
$ws.Range("N30").Value = -49501.6

$ws.Range("H31").Value = 167128.14
$ws.Range("I31").Value = 1418.9688
$ws.Range("J31").Value = 265326.16
$ws.Range("K31").Value = 1418.9688
$ws.Range("L31").Value = 265326.16
$ws.Range("M31").Value = -1123.9688
$ws.Range("N31").Value = -265916.16

$ws.Range("H34").Value = 167128.14
$ws.Range("I34").Value = 1418.9688
$ws.Range("J34").Value = 265326.16
$ws.Range("K34").Value = 1418.9688
$ws.Range("L34").Value = 265326.16
$ws.Range("M34").Value = -1216.9688
$ws.Range("N34").Value = -265730.16

$ws.Range("H112").Value = 29742.834
$ws.Range("J112").Value = 29742.834
$ws.Range("L112").Value = 29742.834
$ws.Range("N112").Value = -32696.834

$ws.Range("H116").Value = 42364.25
$ws.Range("J116").Value = 42364.25
$ws.Range("L116").Value = 42364.25
$ws.Range("N116").Value = -51542.25

$ws.Range("H119").Value = 41250.668
$ws.Range("J119").Value = 41250.668
$ws.Range("L119").Value = 41250.668
$ws.Range("N119").Value = -50926.668

$ws.Range("H128").Value = 49319.6
$ws.Range("J128").Value = 49319.6
$ws.Range("L128").Value = 49319.6
$ws.Range("N128").Value = -59279.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 6115.625
$ws.Range("I3").Value = 3119.3333
$ws.Range("J3").Value = 7913.4
$ws.Range("K3").Value = 9357.999899999999
$ws.Range("L3").Value = 23740.2
$ws.Range("M3").Value = -9245.999899999999
$ws.Range("N3").Value = -23964.2

$ws.Range("H131").Value = 3794.8684
$ws.Range("J131").Value = 1525.4615
$ws.Range("L131").Value = 4576.3845
$ws.Range("N131").Value = -14656.3845

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H96").Value = 23992
$ws.Range("J96").Value = 23992
$ws.Range("L96").Value = 23992
$ws.Range("N96").Value = -29484

$ws.Range("H110").Value = 47702
$ws.Range("J110").Value = 47702
$ws.Range("L110").Value = 47702
$ws.Range("N110").Value = -55882

$ws.Range("H114").Value = 39740.75
$ws.Range("J114").Value = 39740.75
$ws.Range("L114").Value = 39740.75
$ws.Range("N114").Value = -48418.75

$ws.Range("H130").Value = 45940.145
$ws.Range("J130").Value = 45940.145
$ws.Range("L130").Value = 45940.145
$ws.Range("N130").Value = -55980.145

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2743.1875
$ws.Range("I40").Value = 2136.1428
$ws.Range("K40").Value = 2136.1428
$ws.Range("M40").Value = -2000.1428

$ws.Range("H108").Value = 36809
$ws.Range("J108").Value = 36809
$ws.Range("L108").Value = 36809
$ws.Range("N108").Value = -44489

$ws.Range("H110").Value = 22214.666
$ws.Range("J110").Value = 22214.666
$ws.Range("L110").Value = 22214.666
$ws.Range("N110").Value = -30394.666

$ws.Range("H112").Value = 31268.5
$ws.Range("I112").Value = 15000
$ws.Range("K112").Value = 15000
$ws.Range("M112").Value = -13523

$ws.Range("H119").Value = 36206
$ws.Range("J119").Value = 36206
$ws.Range("L119").Value = 36206
$ws.Range("N119").Value = -45882

$ws.Range("H120").Value = 46390
$ws.Range("J120").Value = 46390
$ws.Range("L120").Value = 46390
$ws.Range("N120").Value = -56066

$ws.Range("H124").Value = 47421
$ws.Range("J124").Value = 47421
$ws.Range("L124").Value = 47421
$ws.Range("N124").Value = -57241

$ws.Range("H127").Value = 50702.332
$ws.Range("J127").Value = 50702.332
$ws.Range("L127").Value = 50702.332
$ws.Range("N127").Value = -60622.332

$ws.Range("H128").Value = 32171.6
$ws.Range("J128").Value = 32171.6
$ws.Range("L128").Value = 32171.6
$ws.Range("N128").Value = -42131.6

$ws.Range("H130").Value = 43214.5
$ws.Range("J130").Value = 43214.5
$ws.Range("L130").Value = 43214.5
$ws.Range("N130").Value = -53254.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H108").Value = 32558.666
$ws.Range("J108").Value = 32558.666
$ws.Range("L108").Value = 32558.666
$ws.Range("N108").Value = -40238.666

$ws.Range("H110").Value = 19883
$ws.Range("J110").Value = 19883
$ws.Range("L110").Value = 19883
$ws.Range("N110").Value = -28063

$ws.Range("H116").Value = 49680
$ws.Range("J116").Value = 49680
$ws.Range("L116").Value = 49680
$ws.Range("N116").Value = -58858

$ws.Range("H120").Value = 35210
$ws.Range("J120").Value = 35210
$ws.Range("L120").Value = 35210
$ws.Range("N120").Value = -44886

$ws.Range("H128").Value = 49707
$ws.Range("J128").Value = 49707
$ws.Range("L128").Value = 49707
$ws.Range("N128").Value = -59667

$ws.Range("H131").Value = 47590.4
$ws.Range("J131").Value = 47590.4
$ws.Range("L131").Value = 47590.4
$ws.Range("N131").Value = -57670.4

$ws.Range("H132").Value = 1474
$ws.Range("I132").Value = 848.4375
$ws.Range("J132").Value = 2308.0833
$ws.Range("K132").Value = 2545.3125
$ws.Range("L132").Value = 6924.249899999999
$ws.Range("M132").Value = -15.3125
$ws.Range("N132").Value = -11984.2499
